$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the value in E11 but keep its style (matches Excel's "Delete" key
# behaviour - the cell keeps formatting, just loses its content).
$ws.Range("E11").Value = ""

# Row 12 (the trailing, otherwise-empty row) is removed entirely.
$ws.Rows(12).Delete()

# Update the saved selection/active cell, as recorded in the sheet view.
$ws.Range("D15").Select()
